$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-7 from 45208 to 45212
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45212
}
